$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-11 from serial 45175 to 45177,
# keeping the existing date format/style on the cells.
$ws.Range("C2:C11").Value = 45177
